# Adds 30 new product rows (rows 32-61) to the Products worksheet,
# mirroring additional coffee table listings appended to the report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Lavish Rectangular Metal Coffee Table In Powder Coating Finish"
$ws.Range("C32").Value = "₹23,220"

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Anny Rectangular Metal Coffee Table In Powder Coating Finish"
$ws.Range("C33").Value = "₹23,220"

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Frazer Rectangular Metal Coffee Table In Powder Coating Finish"
$ws.Range("C34").Value = "₹25,451"

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "Marten Rectangular Metal Coffee Table In Powder Coating Finish"
$ws.Range("C35").Value = "₹23,853"

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "Alix Rectangular Metal Coffee Table In Stainless Steel Finish"
$ws.Range("C36").Value = "₹19,763"

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Hazel Round Metal Coffee Table In Powder Coating Finish"
$ws.Range("C37").Value = "₹19,763"

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "Peter Square Metal Coffee Table In Powder Coating Finish"
$ws.Range("C38").Value = "₹19,763"

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "Doug Square Metal Coffee Table In Powder Coating Finish"
$ws.Range("C39").Value = "₹19,763"

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "Angel Round Metal Coffee Table In Powder Coating Finish"
$ws.Range("C40").Value = "₹20,656"

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "Blane Square Solid Wood Coffee Table In Teak Finish"
$ws.Range("C41").Value = "₹19,999"

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Blane Square Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C42").Value = "₹19,499"

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "Blane Square Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C43").Value = "₹19,499"

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "Blane Square Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C44").Value = "₹19,499"

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "Nashville Round Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C45").Value = "₹18,998"

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "Nashville Round Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C46").Value = "₹18,998"

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "Nashville Round Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C47").Value = "₹18,998"

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "Nashville Round Solid Wood Coffee Table In Antique Grey Finish"
$ws.Range("C48").Value = "₹18,998"

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "Blane Square Solid Wood Coffee Table In Teak Finish"
$ws.Range("C49").Value = "₹19,499"

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "Blane Square Solid Wood Coffee Table In Walnut Finish"
$ws.Range("C50").Value = "₹19,499"

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "Nashville Round Solid Wood Coffee Table In Teak Finish"
$ws.Range("C51").Value = "₹18,998"

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "Nashville Round Solid Wood Coffee Table In Teak Finish"
$ws.Range("C52").Value = "₹18,998"

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "Nashville Round Solid Wood Coffee Table In Teak Finish"
$ws.Range("C53").Value = "₹18,998"

$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "Nashville Round Solid Wood Coffee Table In Teak Finish"
$ws.Range("C54").Value = "₹18,998"

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "Irish Rectangular Solid Wood Coffee Table In Walnut Finish"
$ws.Range("C55").Value = "₹13,999"

$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "Irish Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C56").Value = "₹13,999"

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "Irish Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C57").Value = "₹13,999"

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "Irish Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C58").Value = "₹13,999"

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "Irish Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C59").Value = "₹13,999"

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "Milan Square Solid Wood Coffee Table In Walnut Finish"
$ws.Range("C60").Value = "₹13,498"

$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "Montreal Square Solid Wood Coffee Table In Walnut Finish"
$ws.Range("C61").Value = "₹11,498"
